$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("514:515").Insert()

$ws.Cells.Item(514, 1).Value = 9
$ws.Cells.Item(514, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(514, 3).Value = "Metropolitana"
$ws.Cells.Item(514, 4).Value = 44776
$ws.Cells.Item(514, 5).Value = 13
$ws.Cells.Item(514, 6).Value = 100114014
$ws.Cells.Item(514, 7).Value = "Betarraga"
$ws.Cells.Item(514, 8).Value = "Sin especificar"
$ws.Cells.Item(514, 9).Value = "Primera"
$ws.Cells.Item(514, 10).Value = 7000
$ws.Cells.Item(514, 11).Value = 140
$ws.Cells.Item(514, 12).Value = 150
$ws.Cells.Item(514, 13).Value = 145
$ws.Cells.Item(514, 14).Value = '$/unidad'
$ws.Cells.Item(514, 15).Value = "Región Metropolitana"
$ws.Cells.Item(514, 16).Value = 145
$ws.Cells.Item(514, 17).Value = 1
$ws.Cells.Item(514, 18).Value = "Hortaliza"

$ws.Cells.Item(515, 1).Value = 9
$ws.Cells.Item(515, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(515, 3).Value = "Metropolitana"
$ws.Cells.Item(515, 4).Value = 44776
$ws.Cells.Item(515, 5).Value = 13
$ws.Cells.Item(515, 6).Value = 100114014
$ws.Cells.Item(515, 7).Value = "Betarraga"
$ws.Cells.Item(515, 8).Value = "Sin especificar"
$ws.Cells.Item(515, 9).Value = "Segunda"
$ws.Cells.Item(515, 10).Value = 2500
$ws.Cells.Item(515, 11).Value = 120
$ws.Cells.Item(515, 12).Value = 120
$ws.Cells.Item(515, 13).Value = 120
$ws.Cells.Item(515, 14).Value = '$/unidad'
$ws.Cells.Item(515, 15).Value = "Región Metropolitana"
$ws.Cells.Item(515, 16).Value = 120
$ws.Cells.Item(515, 17).Value = 1
$ws.Cells.Item(515, 18).Value = "Hortaliza"
